$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "RTDO L"
$ws.Range("D1").Value = "RTDO V"
$ws.Range("F1").Value = "RTDO L.1"
$ws.Range("I1").Value = "RTDO V.1"
